# Apply updates to the Real-Time streamlit report workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 6
$summary.Range("B4").Value = 10
$summary.Range("B9").Value = 1
$summary.Range("C9").Value = 3
$summary.Range("B10").Value = 5
$summary.Range("C10").Value = 7
$summary.Range("C16").Value = 1
$summary.Range("B17").Value = 1

# --- Sheet: Hourly Breakdown ---
$hourly = $wb.Worksheets.Item("Hourly Breakdown")
$hourly.Range("B16").Value = 0
$hourly.Range("C16").Value = 0
$hourly.Range("D16").Value = 0
$hourly.Range("B17").Value = 0
$hourly.Range("C17").Value = 0
$hourly.Range("B18").Value = 6
$hourly.Range("C18").Value = 10
$hourly.Range("D18").Value = -4

# --- Sheet: Charts Data ---
$charts = $wb.Worksheets.Item("Charts Data")
$charts.Range("B3").Value = 1
$charts.Range("B8").Value = 1
